# yupao p18 24:12 P4
#
# - Re-apply the "Normal" cell style across the existing data range (A1:B5).
#   This flips the style actually used by those cells to a new cellXfs entry
#   (still General number format / default font), matching the workbook's
#   new second cellXfs record.
# - Row 5's member ("zhangsan" in B5) is removed.
# - A new row is appended: member #5, nickname "zhangsan".
# - Selection ends up parked at G13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the Normal style to the current data so it picks up a fresh xf.
$ws.Range("A1:B5").Style = "Normal"

# Remove the nickname that was in B5.
$ws.Range("B5").ClearContents()

# Append a new member row.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "zhangsan"

# Leave the selection where the author left off.
$ws.Range("G13").Select()
